$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 90 to hold a new, most-recent weekly price
# report (shifts the existing rows 90-102 down to 93-105, matching the
# rest of the historical data already present below).
$ws.Rows("90:92").Insert()

# Row 90: Chirimoya "Especial" quality, Provincia de Limarí
$ws.Cells.Item(90, 1).Value = 3
$ws.Cells.Item(90, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(90, 3).Value = "Coquimbo"
$ws.Cells.Item(90, 4).Value = "2021-11-05"
$ws.Cells.Item(90, 5).Value = 5
$ws.Cells.Item(90, 6).Value = "Fruta"
$ws.Cells.Item(90, 7).Value = 100107
$ws.Cells.Item(90, 8).Value = "Otros"
$ws.Cells.Item(90, 9).Value = 100107002
$ws.Cells.Item(90, 10).Value = "Chirimoya"
$ws.Cells.Item(90, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(90, 12).Value = "Especial"
$ws.Cells.Item(90, 13).Value = 45
$ws.Cells.Item(90, 14).Value = 26000
$ws.Cells.Item(90, 15).Value = 26000
$ws.Cells.Item(90, 16).Value = 26000
$ws.Cells.Item(90, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(90, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(90, 19).Value = 2600
$ws.Cells.Item(90, 20).Value = 10

# Row 91: Chirimoya "Primera" quality, Provincia de Limarí
$ws.Cells.Item(91, 1).Value = 3
$ws.Cells.Item(91, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(91, 3).Value = "Coquimbo"
$ws.Cells.Item(91, 4).Value = "2021-11-05"
$ws.Cells.Item(91, 5).Value = 5
$ws.Cells.Item(91, 6).Value = "Fruta"
$ws.Cells.Item(91, 7).Value = 100107
$ws.Cells.Item(91, 8).Value = "Otros"
$ws.Cells.Item(91, 9).Value = 100107002
$ws.Cells.Item(91, 10).Value = "Chirimoya"
$ws.Cells.Item(91, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(91, 12).Value = "Primera"
$ws.Cells.Item(91, 13).Value = 50
$ws.Cells.Item(91, 14).Value = 23000
$ws.Cells.Item(91, 15).Value = 23000
$ws.Cells.Item(91, 16).Value = 23000
$ws.Cells.Item(91, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(91, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(91, 19).Value = 2300
$ws.Cells.Item(91, 20).Value = 10

# Row 92: Chirimoya "Segunda" quality, Provincia de Limarí
$ws.Cells.Item(92, 1).Value = 3
$ws.Cells.Item(92, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(92, 3).Value = "Coquimbo"
$ws.Cells.Item(92, 4).Value = "2021-11-05"
$ws.Cells.Item(92, 5).Value = 5
$ws.Cells.Item(92, 6).Value = "Fruta"
$ws.Cells.Item(92, 7).Value = 100107
$ws.Cells.Item(92, 8).Value = "Otros"
$ws.Cells.Item(92, 9).Value = 100107002
$ws.Cells.Item(92, 10).Value = "Chirimoya"
$ws.Cells.Item(92, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(92, 12).Value = "Segunda"
$ws.Cells.Item(92, 13).Value = 50
$ws.Cells.Item(92, 14).Value = 20000
$ws.Cells.Item(92, 15).Value = 20000
$ws.Cells.Item(92, 16).Value = 20000
$ws.Cells.Item(92, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(92, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(92, 19).Value = 2000
$ws.Cells.Item(92, 20).Value = 10
